# Apply the edits described by the diff:
#  1. Row 11: fill in B11 ("done") and C11 ("0.4") to match the style/format
#     used by the other rows (copy formatting + text type from row 10).
#  2. Row 25: fix typo in the "outsource" task name
#     ("AppVerrior" -> "AppVeyor").
#  3. Add a new row 27 ("Validate yaml" / "done" / "0.4"), duplicated from
#     the existing row 26 so that style/number formatting match exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 11 - add status + version columns
$ws.Range("B10:C10").Copy($ws.Range("B11:C11"))

# 2) New row 27 - "Validate yaml" task, duplicate formatting from row 26
$ws.Range("A26:C26").Copy($ws.Range("A27:C27"))
$ws.Range("A27").Value = "Validate yaml"

# 3) Row 25 - fix typo in task name
$ws.Range("A25").Value = "outsource KsWare.AppVeyor.Api"
